$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "61.647.13"
$ws.Range("E2").Value = "  -2.22%  "
Set-TextValue "D3" "2.949.54"
$ws.Range("E3").Value = "  -3.44%  "
$ws.Range("E4").Value = "  +0.26%  "
Set-TextValue "D5" "582.35"
$ws.Range("E5").Value = "  -0.56%  "
Set-TextValue "D6" "141.82"
$ws.Range("E6").Value = "  -6.60%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -2.85%  "
Set-TextValue "D9" "2.949.92"
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("E10").Value = "  -5.99%  "
$ws.Range("E11").Value = "  -2.50%  "
Set-TextValue "D12" "0.457"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("E13").Value = "  -4.06%  "
Set-TextValue "D14" "33.94"
$ws.Range("E14").Value = "  -6.16%  "
$ws.Range("E15").Value = "  +1.44%  "
Set-TextValue "D16" "3.444.05"
$ws.Range("E16").Value = "  -3.14%  "
Set-TextValue "D17" "6.96"
$ws.Range("E17").Value = "  -2.15%  "
Set-TextValue "D18" "61.645.94"
$ws.Range("E18").Value = "  -2.14%  "
Set-TextValue "D19" "2.949.59"
$ws.Range("E19").Value = "  -3.37%  "
Set-TextValue "D20" "448.47"
$ws.Range("E20").Value = "  -6.19%  "
Set-TextValue "D21" "13.77"
$ws.Range("E21").Value = "  -3.46%  "
Set-TextValue "D22" "0.676"
$ws.Range("E22").Value = "  -4.11%  "
Set-TextValue "D23" "7.24"
$ws.Range("E23").Value = "  -3.41%  "
Set-TextValue "D24" "81.09"
$ws.Range("E24").Value = "  -0.94%  "
Set-TextValue "D25" "12.07"
$ws.Range("E25").Value = "  -4.28%  "
$ws.Range("E26").Value = "  -10.76%  "
$ws.Range("E27").Value = "  -0.04%  "
Set-TextValue "D28" "9.47"
$ws.Range("E28").Value = "  -9.83%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -1.99%  "
Set-TextValue "D31" "6.82"
$ws.Range("E31").Value = "  -7.46%  "
$ws.Range("E32").Value = "  -6.52%  "
Set-TextValue "D33" "27.12"
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("E34").Value = "  -3.87%  "
Set-TextValue "D35" "1.01"
$ws.Range("E35").Value = "  -5.17%  "
Set-TextValue "D36" "0.0₃0769"
$ws.Range("E36").Value = "  -5.70%  "
Set-TextValue "D37" "5.64"
$ws.Range("E37").Value = "  -4.50%  "
Set-TextValue "D38" "2.07"
$ws.Range("E38").Value = "  -6.18%  "
Set-TextValue "D39" "49.96"
$ws.Range("E39").Value = "  -0.86%  "
Set-TextValue "D40" "9.06"
$ws.Range("E40").Value = "  -1.83%  "
Set-TextValue "D41" "0.119"
$ws.Range("E41").Value = "  +3.70%  "
Set-TextValue "D42" "2.78"
$ws.Range("E42").Value = "  -14.35%  "
Set-TextValue "D43" "387.67"
$ws.Range("E43").Value = "  -9.92%  "
$ws.Range("E44").Value = "  -2.79%  "
Set-TextValue "D45" "2.708.52"
$ws.Range("E45").Value = "  -4.30%  "
Set-TextValue "D46" "0.261"
$ws.Range("E46").Value = "  -9.12%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D47" "36.63"
$ws.Range("E47").Value = "  -4.17%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D48" "129.80"
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("E50").Value = "  -1.74%  "
Set-TextValue "D51" "2.15"
$ws.Range("E51").Value = "  -2.00%  "
